# aula do dia 26/03
# Adds the 26/03 (and surrounding 25/03) class columns to the FREQ
# attendance sheet: a new 5-column block (PDMO, RMST, PROJ, PROJ, PROJ,
# PDMO, RMST) mirroring the existing weekly pattern, with the 25/03 and
# 26/03 dates filled in on row 2 and the first two attendance marks
# (column AA) filled in for each student.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FREQ")

# --- Row 1: header labels for the new columns (AA..AG) ---------------
$headerVals = @("PDMO", "RMST", "PROJ", "PROJ", "PROJ", "PDMO", "RMST")
$col = 27   # column AA
foreach ($v in $headerVals) {
    $ws.Cells.Item(1, $col).Value = $v
    $col = $col + 1
}

# --- Row 2: class dates. Only the first two (25/03, 26/03) are known -
# the remaining five columns of the block are pre-formatted (date
# number format) but left blank, same as the existing weekly blocks
# before they get filled in.
for ($c = 27; $c -le 33; $c++) {
    $ws.Cells.Item(2, $c).NumberFormat = "d-mmm"
}
$ws.Cells.Item(2, 27).Value = 44280   # AA2 = 25-Mar-2021
$ws.Cells.Item(2, 28).Value = 44281   # AB2 = 26-Mar-2021

# --- Column AA: attendance mark for the 25/03 class, per student -----
# Row 5 (CRISTIAN RAFAEL DA SILVA FERREIRA) has no mark for this class
# yet, so it is intentionally left untouched.
$attendance = @{
    3  = "F"
    4  = "P"
    6  = "P"
    7  = "F"
    8  = "P"
    9  = "P"
    10 = "p"
    11 = "P"
    12 = "F"
    13 = "P"
    14 = "P"
    15 = "P"
    16 = "P"
    17 = "P"
    18 = "P"
    19 = "P"
    20 = "P"
}

foreach ($r in $attendance.Keys) {
    $ws.Cells.Item($r, 27).Value = $attendance[$r]
}

# --- Leave the selection where the author left it when saving --------
$ws.Range("AA11").Select() | Out-Null
